$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13/14: swap the INTERCEPT / SLOPE argument order, and update B14 ---
$ws.Range("D13").Formula = "=INTERCEPT(B13:B14,A13:A14)"
$ws.Range("B14").Value = 15086
$ws.Range("D14").Formula = "=SLOPE(B13:B14,A13:A14)"

# --- Remove old row 17 (C17 = 8746*0.612 + 6) ---
$ws.Range("C17").ClearContents()

# --- New rows 15/16: INTERCEPT/SLOPE over empty ranges -> #DIV/0! ---
$ws.Range("D15").Formula = "=INTERCEPT(B15:B16,A15:A16)"
$ws.Range("D16").Formula = "=SLOPE(B15:B16,A15:A16)"

# --- New rows 18/19 ---
$ws.Range("A18").Value = 141
$ws.Range("B18").Value = 225
$ws.Range("D18").Formula = "=INTERCEPT(B18:B19,A18:A19)"
$ws.Range("G18").Formula = "=1/0.055"

$ws.Range("A19").Value = 9367
$ws.Range("B19").Value = 18200
$ws.Range("D19").Formula = "=SLOPE(B18:B19,A18:A19)"
$ws.Range("G19").Formula = "=1/4.45"

# --- New rows 21/22 ---
$ws.Range("A21").Value = 68
$ws.Range("B21").Value = 25
$ws.Range("D21").Formula = "=INTERCEPT(B21:B22,A21:A22)"

$ws.Range("A22").Value = 2518
$ws.Range("B22").Value = 1000
$ws.Range("D22").Formula = "=SLOPE(B21:B22,A21:A22)"

# --- Update the active selection shown when the sheet was last saved ---
$ws.Range("E19").Select()
